# GPLIM-4825 ignore the pipeline for now and clean up DB model
# Append a new data row (Eppendorf96 / 77891 / UMI 2 / Spacer 2 / Before First Read)
# to the bottom of the UMI reagents table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Eppendorf96"
$ws.Range("B8").Value = 77891
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Before First Read"

$ws.Range("E8").Select()
